$d = $word.ActiveDocument

$replacements = @(
    @("98×34=", "65×63="),
    @("96×40=", "93×78="),
    @("18×56=", "60×27="),
    @("13×50=", "62×18="),
    @("49×23=", "64×11="),
    @("66×25=", "54×18="),
    @("55×15=", "79×32="),
    @("41×37=", "22×29="),
    @("53×13=", "67×33="),
    @("68×51=", "27×65="),
    @("97×12=", "70×24="),
    @("87×71=", "32×81="),
    @("36×91=", "70×90="),
    @("77×73=", "38×65="),
    @("77×97=", "34×99="),
    @("32×94=", "11×15="),
    @("31×69=", "39×76="),
    @("42×42=", "97×53="),
    @("15×50=", "15×27="),
    @("54×51=", "69×67="),
    @("48×78=", "43×12="),
    @("87×57=", "33×44="),
    @("68×89=", "40×71="),
    @("24×73=", "64×55="),
    @("69×28=", "72×26=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
